# Insert a new data row at row 186 (pushing the existing rows 186:281 down
# to 187:282) and populate it with a new price-report record for Agrícola
# del Norte S.A. de Arica - Limón (Tahití, Primera).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 186:281 down one position, creating a blank row 186.
$ws.Rows("186:186").Insert()

# Fill the newly inserted row with the new record.
$ws.Range("A186").Value = 1
$ws.Range("B186").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C186").Value = "Arica y Parinacota"
$ws.Range("D186").Value = 44813
$ws.Range("E186").Value = 15
$ws.Range("F186").Value = "Fruta"
$ws.Range("G186").Value = 100102
$ws.Range("H186").Value = "Cítricos"
$ws.Range("I186").Value = 100102003
$ws.Range("J186").Value = "Limón"
$ws.Range("K186").Value = "Tahití"
$ws.Range("L186").Value = "Primera"
$ws.Range("M186").Value = 144
$ws.Range("N186").Value = 53000
$ws.Range("O186").Value = 54000
$ws.Range("P186").Value = 53500
$ws.Range("Q186").Value = "$/caja 24 kilos"
$ws.Range("R186").Value = "Perú"
$ws.Range("S186").Value = 2229
$ws.Range("T186").Value = 24
